$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecasts")

$ws.Range("B4").Value = 0.9973508863433136
$ws.Range("E4").Value = 0.01918230390331595

$ws.Range("B5").Value = 68.62472964229038
$ws.Range("E5").Value = 1.227221335282634

$ws.Range("B6").Value = 0.4697089
$ws.Range("C6").Value = 0.5302911
$ws.Range("E6").Value = 0.9981807
$ws.Range("F6").Value = 0.0018193

$ws.Range("B7").Value = 26.5056628
$ws.Range("C7").Value = 26.944527
$ws.Range("E7").Value = 32.5121498
$ws.Range("F7").Value = 16.0645064
